$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly added values for row 5 and row 6 (columns A and B)
$ws.Range("A5").Value = 9321944.0
$ws.Range("B5").Value = 7990296.0

$ws.Range("A6").Value = 7874784.0
$ws.Range("B6").Value = 9301696.0

# Add new rows 9, 10, 11 with values for columns C and D
$ws.Range("C9").Value = 1331800.0
$ws.Range("D9").Value = 1331712.0

$ws.Range("C10").Value = 1331688.0
$ws.Range("D10").Value = 1331688.0

$ws.Range("C11").Value = 1331744.0
$ws.Range("D11").Value = 1331688.0
